$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 9
$ws.Range("B5").Value = 9
$ws.Range("B6").Value = 9
$ws.Range("B7").Value = 9

$ws.Range("C9").Select()
